# Natmi following Dr Hou advice
# Recompute the Fgf2 -> Fgfr3 ligand/receptor edge table with updated
# ligand/receptor-expressing-cell counts (now 3 cells per cluster instead
# of 1) and the resulting recalculated expression/specificity metrics.
# This also adds three new target-cluster rows for the "sCs" sending
# cluster (rows 8-10), growing the table from A1:T7 to A1:T10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fgf2"
$ws.Cells.Item(2, 3).Value = "Fgfr3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.747119
$ws.Cells.Item(2, 8).Value = 2.241357
$ws.Cells.Item(2, 9).Value = 0.03096954854571248
$ws.Cells.Item(2, 10).Value = 0.03096954854571248
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.656156333333334
$ws.Cells.Item(2, 14).Value = 4.968469000000001
$ws.Cells.Item(2, 15).Value = 0.6151212440816572
$ws.Cells.Item(2, 16).Value = 0.6151212440816572
$ws.Cells.Item(2, 17).Value = 1.237345863603667
$ws.Cells.Item(2, 18).Value = 11.136112772433
$ws.Cells.Item(2, 19).Value = 0.01905002723008594
$ws.Cells.Item(2, 20).Value = 0.01905002723008594
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fgf2"
$ws.Cells.Item(3, 3).Value = "Fgfr3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.747119
$ws.Cells.Item(3, 8).Value = 2.241357
$ws.Cells.Item(3, 9).Value = 0.03096954854571248
$ws.Cells.Item(3, 10).Value = 0.03096954854571248
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.5648773333333333
$ws.Cells.Item(3, 14).Value = 1.694632
$ws.Cells.Item(3, 15).Value = 0.2098038941373262
$ws.Cells.Item(3, 16).Value = 0.2098038941373262
$ws.Cells.Item(3, 17).Value = 0.4220305884026667
$ws.Cells.Item(3, 18).Value = 3.798275295623999
$ws.Cells.Item(3, 19).Value = 0.006497531884565445
$ws.Cells.Item(3, 20).Value = 0.006497531884565444
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fgf2"
$ws.Cells.Item(4, 3).Value = "Fgfr3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.747119
$ws.Cells.Item(4, 8).Value = 2.241357
$ws.Cells.Item(4, 9).Value = 0.03096954854571248
$ws.Cells.Item(4, 10).Value = 0.03096954854571248
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.4713726666666667
$ws.Cells.Item(4, 14).Value = 1.414118
$ws.Cells.Item(4, 15).Value = 0.1750748617810164
$ws.Cells.Item(4, 16).Value = 0.1750748617810165
$ws.Cells.Item(4, 17).Value = 0.3521714753473333
$ws.Cells.Item(4, 18).Value = 3.169543278126
$ws.Cells.Item(4, 19).Value = 0.005421989431061091
$ws.Cells.Item(4, 20).Value = 0.005421989431061092
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fgf2"
$ws.Cells.Item(5, 3).Value = "Fgfr3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 19.74619233333334
$ws.Cells.Item(5, 8).Value = 59.23857700000001
$ws.Cells.Item(5, 9).Value = 0.8185184181638298
$ws.Cells.Item(5, 10).Value = 0.8185184181638298
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.656156333333334
$ws.Cells.Item(5, 14).Value = 4.968469000000001
$ws.Cells.Item(5, 15).Value = 0.6151212440816572
$ws.Cells.Item(5, 16).Value = 0.6151212440816572
$ws.Cells.Item(5, 17).Value = 32.70278149206812
$ws.Cells.Item(5, 18).Value = 294.3250334286131
$ws.Cells.Item(5, 19).Value = 0.5034880676846851
$ws.Cells.Item(5, 20).Value = 0.5034880676846851
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fgf2"
$ws.Cells.Item(6, 3).Value = "Fgfr3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 19.74619233333334
$ws.Cells.Item(6, 8).Value = 59.23857700000001
$ws.Cells.Item(6, 9).Value = 0.8185184181638298
$ws.Cells.Item(6, 10).Value = 0.8185184181638298
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.5648773333333333
$ws.Cells.Item(6, 14).Value = 1.694632
$ws.Cells.Item(6, 15).Value = 0.2098038941373262
$ws.Cells.Item(6, 16).Value = 0.2098038941373262
$ws.Cells.Item(6, 17).Value = 11.15417646874045
$ws.Cells.Item(6, 18).Value = 100.387588218664
$ws.Cells.Item(6, 19).Value = 0.1717283515538958
$ws.Cells.Item(6, 20).Value = 0.1717283515538958
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fgf2"
$ws.Cells.Item(7, 3).Value = "Fgfr3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 19.74619233333334
$ws.Cells.Item(7, 8).Value = 59.23857700000001
$ws.Cells.Item(7, 9).Value = 0.8185184181638298
$ws.Cells.Item(7, 10).Value = 0.8185184181638298
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.4713726666666667
$ws.Cells.Item(7, 14).Value = 1.414118
$ws.Cells.Item(7, 15).Value = 0.1750748617810164
$ws.Cells.Item(7, 16).Value = 0.1750748617810165
$ws.Cells.Item(7, 17).Value = 9.307815336676224
$ws.Cells.Item(7, 18).Value = 83.77033803008602
$ws.Cells.Item(7, 19).Value = 0.1433019989252487
$ws.Cells.Item(7, 20).Value = 0.1433019989252487
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Fgf2"
$ws.Cells.Item(8, 3).Value = "Fgfr3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 3.630999
$ws.Cells.Item(8, 8).Value = 10.892997
$ws.Cells.Item(8, 9).Value = 0.1505120332904577
$ws.Cells.Item(8, 10).Value = 0.1505120332904577
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.656156333333334
$ws.Cells.Item(8, 14).Value = 4.968469000000001
$ws.Cells.Item(8, 15).Value = 0.6151212440816572
$ws.Cells.Item(8, 16).Value = 0.6151212440816572
$ws.Cells.Item(8, 17).Value = 6.013501990177001
$ws.Cells.Item(8, 18).Value = 54.121517911593
$ws.Cells.Item(8, 19).Value = 0.09258314916688613
$ws.Cells.Item(8, 20).Value = 0.09258314916688615
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Fgf2"
$ws.Cells.Item(9, 3).Value = "Fgfr3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 3.630999
$ws.Cells.Item(9, 8).Value = 10.892997
$ws.Cells.Item(9, 9).Value = 0.1505120332904577
$ws.Cells.Item(9, 10).Value = 0.1505120332904577
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.5648773333333333
$ws.Cells.Item(9, 14).Value = 1.694632
$ws.Cells.Item(9, 15).Value = 0.2098038941373262
$ws.Cells.Item(9, 16).Value = 0.2098038941373262
$ws.Cells.Item(9, 17).Value = 2.051069032456
$ws.Cells.Item(9, 18).Value = 18.459621292104
$ws.Cells.Item(9, 19).Value = 0.0315780106988649
$ws.Cells.Item(9, 20).Value = 0.0315780106988649
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Fgf2"
$ws.Cells.Item(10, 3).Value = "Fgfr3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 3.630999
$ws.Cells.Item(10, 8).Value = 10.892997
$ws.Cells.Item(10, 9).Value = 0.1505120332904577
$ws.Cells.Item(10, 10).Value = 0.1505120332904577
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.4713726666666667
$ws.Cells.Item(10, 14).Value = 1.414118
$ws.Cells.Item(10, 15).Value = 0.1750748617810164
$ws.Cells.Item(10, 16).Value = 0.1750748617810165
$ws.Cells.Item(10, 17).Value = 1.711553681294
$ws.Cells.Item(10, 18).Value = 15.403983131646
$ws.Cells.Item(10, 19).Value = 0.02635087342470662
$ws.Cells.Item(10, 20).Value = 0.02635087342470663